$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (row 1) entirely; remaining rows shift up by one.
$ws.Rows.Item(1).Delete()

# Restore the selection to the (now) first row, matching target sqref.
$ws.Range("A1:XFD1").Select()
